{"js": "// The document contains a single 20x5 table of arithmetic expressions\n// (e.g. \"77+10=\"), one per cell, in row-major order. This script replaces\n// each cell's expression with its new value while preserving the existing\n// run/paragraph formatting (font, size, justification).\nconst table = context.document.body.tables.getFirst();\nconst replacements = [[\"77+10=\", \"11-8=\"], [\"98-5=\", \"11-6=\"], [\"78-76=\", \"82-2=\"], [\"89-57=\", \"91-55=\"], [\"93+4=\", \"39-22=\"], [\"46+52=\", \"51+14=\"], [\"77+21=\", \"45+22=\"], [\"95-35=\", \"71-38=\"], [\"45-22=\", \"32-6=\"], [\"5+16=\", \"77-18=\"], [\"4+84=\", \"61-41=\"], [\"81-3=\", \"82-2=\"], [\"40+37=\", \"33-0=\"], [\"60-10=\", \"91-77=\"], [\"57-1=\", \"36-26=\"], [\"63+3=\", \"88-21=\"], [\"27-13=\", \"47+31=\"], [\"71-61=\", \"85-2=\"], [\"85-60=\", \"38+16=\"], [\"12+9=\", \"8+6=\"], [\"8+88=\", \"52-49=\"], [\"3+13=\", \"3+58=\"], [\"90-54=\", \"95-10=\"], [\"83-66=\", \"12+55=\"], [\"67+1=\", \"73-6=\"], [\"16+2=\", \"75+15=\"], [\"4+12=\", \"1+73=\"], [\"88-27=\", \"62-52=\"], [\"88-48=\", \"3+20=\"], [\"64+23=\", \"85+8=\"], [\"34-28=\", \"69+18=\"], [\"69+16=\", \"49-15=\"], [\"26+7=\", \"24+47=\"], [\"66+26=\", \"2+47=\"], [\"81-18=\", \"66-0=\"], [\"69+12=\", \"85-64=\"], [\"36-0=\", \"21+71=\"], [\"26+31=\", \"83-77=\"], [\"3+15=\", \"39+58=\"], [\"58-36=\", \"55+44=\"], [\"5+71=\", \"48+41=\"], [\"36+33=\", \"1+48=\"], [\"1+35=\", \"10+45=\"], [\"46+48=\", \"30+58=\"], [\"47-7=\", \"85-46=\"], [\"83-52=\", \"88-83=\"], [\"59+35=\", \"96-24=\"], [\"14+51=\", \"31-9=\"], [\"25-18=\", \"45-20=\"], [\"53-13=\", \"35+41=\"], [\"79-64=\", \"23+53=\"], [\"17+63=\", \"57-14=\"], [\"54+32=\", \"92+5=\"], [\"89-84=\", \"58+35=\"], [\"60-22=\", \"68+25=\"], [\"34-24=\", \"18+0=\"], [\"0+77=\", \"76+13=\"], [\"10+34=\", \"77-58=\"], [\"49-44=\", \"90-32=\"], [\"76-65=\", \"68-16=\"], [\"17+69=\", \"68-7=\"], [\"18+59=\", \"34+5=\"], [\"21+44=\", \"35-26=\"], [\"11+13=\", \"21+26=\"], [\"75-49=\", \"51+24=\"], [\"88+0=\", \"55-39=\"], [\"54+15=\", \"70+22=\"], [\"66-19=\", \"85-27=\"], [\"17+15=\", \"29-22=\"], [\"65+6=\", \"52+11=\"], [\"70+27=\", \"1+4=\"], [\"94-45=\", \"83-49=\"], [\"25-2=\", \"1+39=\"], [\"31-8=\", \"35-5=\"], [\"33+28=\", \"71+25=\"], [\"17+72=\", \"42-33=\"], [\"91-82=\", \"68-22=\"], [\"93-69=\", \"45+13=\"], [\"74-59=\", \"63+4=\"], [\"56+13=\", \"48-30=\"], [\"23+57=\", \"36+24=\"], [\"57-23=\", \"10+24=\"], [\"20-10=\", \"15+77=\"], [\"19+35=\", \"36+1=\"], [\"75-1=\", \"13+75=\"], [\"37+57=\", \"59+26=\"], [\"62-45=\", \"54+39=\"], [\"61+34=\", \"13-2=\"], [\"98-1=\", \"5+1=\"], [\"75-70=\", \"55-20=\"], [\"10+23=\", \"44+53=\"], [\"99-0=\", \"86-35=\"], [\"18+37=\", \"24+75=\"], [\"28-21=\", \"21+67=\"], [\"7+86=\", \"28+21=\"], [\"75+16=\", \"17+64=\"], [\"67-13=\", \"1+9=\"], [\"35+64=\", \"25-20=\"], [\"42+35=\", \"89-25=\"], [\"30-10=\", \"38+28=\"]];\nconst COLS = 5;\n\n// Grab each target cell's first paragraph up front.\nconst paras = [];\nfor (let i = 0; i < replacements.length; i++) {\n  const row = Math.floor(i / COLS);\n  const col = i % COLS;\n  const cell = table.getCell(row, col);\n  const para = cell.body.paragraphs.getFirst();\n  para.load(\"text\");\n  paras.push(para);\n}\nawait context.sync();\n\n// Replace the text of each paragraph's range with the new expression,\n// using insertText(..., replace) on the paragraph range (not the cell\n// body) so the run's rFonts/sz formatting is preserved.\nfor (let i = 0; i < replacements.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = paras[i];\n  const range = para.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The document contains a single 20x5 table of arithmetic expressions\n# (e.g. \"77+10=\"), one per cell, in row-major order. This script replaces\n# each cell's expression with its new value while preserving the existing\n# run/paragraph formatting (font, size, justification).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$cols = 5\n$oldValues = @('77+10=', '98-5=', '78-76=', '89-57=', '93+4=', '46+52=', '77+21=', '95-35=', '45-22=', '5+16=', '4+84=', '81-3=', '40+37=', '60-10=', '57-1=', '63+3=', '27-13=', '71-61=', '85-60=', '12+9=', '8+88=', '3+13=', '90-54=', '83-66=', '67+1=', '16+2=', '4+12=', '88-27=', '88-48=', '64+23=', '34-28=', '69+16=', '26+7=', '66+26=', '81-18=', '69+12=', '36-0=', '26+31=', '3+15=', '58-36=', '5+71=', '36+33=', '1+35=', '46+48=', '47-7=', '83-52=', '59+35=', '14+51=', '25-18=', '53-13=', '79-64=', '17+63=', '54+32=', '89-84=', '60-22=', '34-24=', '0+77=', '10+34=', '49-44=', '76-65=', '17+69=', '18+59=', '21+44=', '11+13=', '75-49=', '88+0=', '54+15=', '66-19=', '17+15=', '65+6=', '70+27=', '94-45=', '25-2=', '31-8=', '33+28=', '17+72=', '91-82=', '93-69=', '74-59=', '56+13=', '23+57=', '57-23=', '20-10=', '19+35=', '75-1=', '37+57=', '62-45=', '61+34=', '98-1=', '75-70=', '10+23=', '99-0=', '18+37=', '28-21=', '7+86=', '75+16=', '67-13=', '35+64=', '42+35=', '30-10=')\n$newValues = @('11-8=', '11-6=', '82-2=', '91-55=', '39-22=', '51+14=', '45+22=', '71-38=', '32-6=', '77-18=', '61-41=', '82-2=', '33-0=', '91-77=', '36-26=', '88-21=', '47+31=', '85-2=', '38+16=', '8+6=', '52-49=', '3+58=', '95-10=', '12+55=', '73-6=', '75+15=', '1+73=', '62-52=', '3+20=', '85+8=', '69+18=', '49-15=', '24+47=', '2+47=', '66-0=', '85-64=', '21+71=', '83-77=', '39+58=', '55+44=', '48+41=', '1+48=', '10+45=', '30+58=', '85-46=', '88-83=', '96-24=', '31-9=', '45-20=', '35+41=', '23+53=', '57-14=', '92+5=', '58+35=', '68+25=', '18+0=', '76+13=', '77-58=', '90-32=', '68-16=', '68-7=', '34+5=', '35-26=', '21+26=', '51+24=', '55-39=', '70+22=', '85-27=', '29-22=', '52+11=', '1+4=', '83-49=', '1+39=', '35-5=', '71+25=', '42-33=', '68-22=', '45+13=', '63+4=', '48-30=', '36+24=', '10+24=', '15+77=', '36+1=', '13+75=', '59+26=', '54+39=', '13-2=', '5+1=', '55-20=', '44+53=', '86-35=', '24+75=', '21+67=', '28+21=', '17+64=', '1+9=', '25-20=', '89-25=', '38+28=')\nfor ($i = 0; $i -lt $newValues.Count; $i++) {\n    $row = [int][Math]::Floor($i / $cols) + 1\n    $col = ($i % $cols) + 1\n    $cell = $t.Cell($row, $col)\n    # Cell.Range.Text setter replaces the cell content (preserving the\n    # trailing end-of-cell marker) and keeps the existing run formatting.\n    $cell.Range.Text = $newValues[$i]\n}\n"}
